$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: nfo filename changed to the new example path
$ws.Range("A2").Value = "example/example.nfo"

# E2: premiered date is now stored as literal text "2023-05-01"
# (quote-prefixed) instead of a real date serial number.
$e2 = $ws.Range("E2")
$e2.Style = "Normal"
$e2.Value = "'2023-05-01"
$e2.NumberFormat = "mm-dd-yy"

# New column H: movie/director:string -> Spielberg
$ws.Range("H1").Value = "movie/director:string"
$ws.Range("H2").Value = "Spielberg"

# Update the active selection shown when the workbook is reopened
$null = $ws.Range("D6").Select()
